$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,28)
foreach ($r in $rows) {
    $ws.Range("A$r").Value = "x"
}

$ws.Range("B27").Select()
